$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data table (A1:E11) that feeds the line chart.
$ws.Range("E6").Value = 0

$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 0

$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 5

$ws.Range("B9").Value = 19
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 16
$ws.Range("E9").Value = 167

$ws.Range("B10").Value = 39
$ws.Range("C10").Value = 31
$ws.Range("D10").Value = 22
$ws.Range("E10").Value = 234

$ws.Range("B11").Value = 17
$ws.Range("C11").Value = 17
$ws.Range("D11").Value = 27
$ws.Range("E11").Value = 0

# Update the active selection to match the saved view state.
$null = $ws.Range("C15").Select()
